$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column indices: A=1 (rank), B=2 (Coin), C=3 (Link), D=4 (Price), E=5 (Volume(1h))

# Pre-format cells whose new Price value would otherwise be auto-coerced to a
# number by Excel (plain decimal-looking strings) so they remain stored as text,
# matching the source data which always stores Price/Volume as text.
$textFormatCells = @("D5","D6","D8","D9","D10","D11","D12","D13","D15","D16","D19","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D38","D41","D42","D45","D47","D51")
foreach ($ref in $textFormatCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Cells.Item(2, 4).Value = "71.478.37"
$ws.Cells.Item(2, 5).Value = "  +3.10%  "

# Row 3 - Ethereum
$ws.Cells.Item(3, 4).Value = "3.710.72"
$ws.Cells.Item(3, 5).Value = "  +8.01%  "

# Row 4 - TetherUSD
$ws.Cells.Item(4, 5).Value = "  -0.08%  "

# Row 5 - BNB
$ws.Cells.Item(5, 4).Value = "584.37"
$ws.Cells.Item(5, 5).Value = "  +0.67%  "

# Row 6 - Solana
$ws.Cells.Item(6, 4).Value = "179.93"
$ws.Cells.Item(6, 5).Value = "  +1.60%  "

# Row 7 - LidoStakedEther
$ws.Cells.Item(7, 4).Value = "3.701.52"
$ws.Cells.Item(7, 5).Value = "  +7.99%  "

# Row 8 - XRP
$ws.Cells.Item(8, 4).Value = "0.616"
$ws.Cells.Item(8, 5).Value = "  +3.90%  "

# Row 9 - USDC
$ws.Cells.Item(9, 4).Value = "0.999"
$ws.Cells.Item(9, 5).Value = "  -0.12%  "

# Row 10 - Dogecoin
$ws.Cells.Item(10, 4).Value = "0.202"
$ws.Cells.Item(10, 5).Value = "  +2.09%  "

# Row 11 - Cardano
$ws.Cells.Item(11, 4).Value = "0.612"
$ws.Cells.Item(11, 5).Value = "  +4.47%  "

# Row 12 - Avalanche
$ws.Cells.Item(12, 4).Value = "49.58"
$ws.Cells.Item(12, 5).Value = "  +1.79%  "

# Row 13 - ShibaInu
$ws.Cells.Item(13, 4).Value = "0.0000288"
$ws.Cells.Item(13, 5).Value = "  +2.55%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Cells.Item(14, 4).Value = "4.291.98"
$ws.Cells.Item(14, 5).Value = "  +8.40%  "

# Row 15 - BitcoinCash
$ws.Cells.Item(15, 4).Value = "680.53"
$ws.Cells.Item(15, 5).Value = "  -2.80%  "

# Row 16 - Polkadot
$ws.Cells.Item(16, 4).Value = "9.08"
$ws.Cells.Item(16, 5).Value = "  +4.94%  "

# Row 17 - WrappedEther
$ws.Cells.Item(17, 4).Value = "3.704.12"
$ws.Cells.Item(17, 5).Value = "  +8.21%  "

# Row 18 - WrappedBTC
$ws.Cells.Item(18, 4).Value = "71.564.65"
$ws.Cells.Item(18, 5).Value = "  +3.07%  "

# Row 19 - TRON
$ws.Cells.Item(19, 4).Value = "0.123"
$ws.Cells.Item(19, 5).Value = "  +1.20%  "

# Row 20 - Chainlink
$ws.Cells.Item(20, 4).Value = "18.06"
$ws.Cells.Item(20, 5).Value = "  +1.72%  "

# Row 21 - Uniswap
$ws.Cells.Item(21, 4).Value = "11.65"

# Row 22 - Toncoin
$ws.Cells.Item(22, 4).Value = "6.46"
$ws.Cells.Item(22, 5).Value = "  +19.66%  "

# Row 23 - Polygon
$ws.Cells.Item(23, 4).Value = "0.946"
$ws.Cells.Item(23, 5).Value = "  +4.90%  "

# Row 24 - InternetComputer(DFINITY)
$ws.Cells.Item(24, 4).Value = "17.49"
$ws.Cells.Item(24, 5).Value = "  +2.87%  "

# Row 25 - Litecoin
$ws.Cells.Item(25, 4).Value = "102.71"
$ws.Cells.Item(25, 5).Value = "  +1.37%  "

# Row 26 - PancakeSwap
$ws.Cells.Item(26, 4).Value = "4.01"
$ws.Cells.Item(26, 5).Value = "  +2.94%  "

# Row 27 - ImmutableX
$ws.Cells.Item(27, 4).Value = "2.85"
$ws.Cells.Item(27, 5).Value = "  +6.39%  "

# Row 28 - RenderToken
$ws.Cells.Item(28, 4).Value = "10.45"
$ws.Cells.Item(28, 5).Value = "  +8.51%  "

# Row 29 - EthereumClassic
$ws.Cells.Item(29, 4).Value = "35.67"
$ws.Cells.Item(29, 5).Value = "  +5.96%  "

# Row 30 - Filecoin
$ws.Cells.Item(30, 4).Value = "9.21"
$ws.Cells.Item(30, 5).Value = "  +4.93%  "

# Row 31 - NEARProtocol
$ws.Cells.Item(31, 4).Value = "7.40"
$ws.Cells.Item(31, 5).Value = "  +5.58%  "

# Row 32 - dogwifhat
$ws.Cells.Item(32, 4).Value = "4.12"
$ws.Cells.Item(32, 5).Value = "  +9.40%  "

# Row 33 - Bittensor
$ws.Cells.Item(33, 4).Value = "588.06"
$ws.Cells.Item(33, 5).Value = "  +3.28%  "

# Row 34 - Cosmos
$ws.Cells.Item(34, 4).Value = "11.25"
$ws.Cells.Item(34, 5).Value = "  +1.98%  "

# Row 35 - Hedera
$ws.Cells.Item(35, 4).Value = "0.109"
$ws.Cells.Item(35, 5).Value = "  +3.51%  "

# Row 36 - OKB
$ws.Cells.Item(36, 4).Value = "59.21"
$ws.Cells.Item(36, 5).Value = "  +1.80%  "

# Row 37 - Dai
$ws.Cells.Item(37, 5).Value = "  +0.05%  "

# Row 38 - Maker->Kaspa
$ws.Cells.Item(38, 2).Value = "Kaspa"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(38, 4).Value = "0.147"
$ws.Cells.Item(38, 5).Value = "  +5.95%  "

# Row 39 - Kaspa->Maker
$ws.Cells.Item(39, 2).Value = "Maker"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(39, 4).Value = "3.686.03"
$ws.Cells.Item(39, 5).Value = "  +2.18%  "

# Row 40 - PEPE
$ws.Cells.Item(40, 4).Value = "0.0₃0779"
$ws.Cells.Item(40, 5).Value = "  +6.24%  "

# Row 41 - InjectiveProtocol
$ws.Cells.Item(41, 4).Value = "35.73"
$ws.Cells.Item(41, 5).Value = "  +2.29%  "

# Row 42 - Stacks
$ws.Cells.Item(42, 4).Value = "3.47"
$ws.Cells.Item(42, 5).Value = "  +5.94%  "

# Row 43 - Fetch.AI
$ws.Cells.Item(43, 5).Value = "  +4.69%  "

# Row 44 - VeChain
$ws.Cells.Item(44, 5).Value = "  +9.26%  "

# Row 45 - TheGraph
$ws.Cells.Item(45, 4).Value = "0.350"
$ws.Cells.Item(45, 5).Value = "  +4.83%  "

# Row 46 - ThetaToken
$ws.Cells.Item(46, 5).Value = "  +8.34%  "

# Row 47 - ApeXProtocol
$ws.Cells.Item(47, 4).Value = "3.39"
$ws.Cells.Item(47, 5).Value = "  +1.60%  "

# Row 48 - Stellar
$ws.Cells.Item(48, 5).Value = "  +4.24%  "

# Row 49 - Mantle
$ws.Cells.Item(49, 5).Value = "  -0.90%  "

# Row 50 - FirstDigitalUSD
$ws.Cells.Item(50, 5).Value = "  -0.13%  "

# Row 51 - Monero
$ws.Cells.Item(51, 4).Value = "136.20"
$ws.Cells.Item(51, 5).Value = "  +3.70%  "
